# feat: add 2022-Q1 data
#
# The workbook currently has two sheets: "2021-Q4" (sheet1) and "总计"
# (summary sheet). This script inserts a new "2022-Q1" sheet (fund holding
# data, same layout as "2021-Q4") between them, and updates the "总计"
# sheet with a new first data row summarizing "2022-Q1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0. Drop the old "总计" sheet first (it gets rebuilt below, after the
#    new "2022-Q1" sheet) - doing the delete up front keeps sheetId
#    allocation in the same order the source workbook was authored in:
#    2021-Q4=1, 2022-Q1=2, 总计=3.
# ---------------------------------------------------------------------
$wsOldTotal = $wb.Worksheets.Item("总计")
$wsOldTotal.Delete()

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet right after "2021-Q4".
# ---------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item(1)
$wsQ1 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ1.Name = "2022-Q1"

# Clone header-row (B1:H1) and index-column (A2:A8) formatting from the
# "2021-Q4" sheet so the new sheet matches its look (bold/centered/boxed
# style).
$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)

$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A8").PasteSpecial(-4122)

# Header row.
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Column B (fund codes, e.g. "014207") and columns D:G (numeric-looking
# percentages/sizes) are all stored as text, not real numbers - force
# text storage so leading zeros / original formatting survive.
$wsQ1.Range("B2:B8").NumberFormat = "@"
$wsQ1.Range("D2:G8").NumberFormat = "@"

# Index column (A) - sequential 0-based row numbers.
$wsQ1.Range("A2").Value = 0
$wsQ1.Range("A3").Value = 1
$wsQ1.Range("A4").Value = 2
$wsQ1.Range("A5").Value = 3
$wsQ1.Range("A6").Value = 4
$wsQ1.Range("A7").Value = 5
$wsQ1.Range("A8").Value = 6

# Fund rows.
$wsQ1.Range("B2").Value = "014207"
$wsQ1.Range("C2").Value = "华安产业精选混合A"
$wsQ1.Range("D2").Value = "27.31"
$wsQ1.Range("E2").Value = "62.03"
$wsQ1.Range("F2").Value = "1.74"
$wsQ1.Range("G2").Value = "0.4752"
$wsQ1.Range("H2").Value = 9

$wsQ1.Range("B3").Value = "014208"
$wsQ1.Range("C3").Value = "华安产业精选混合C"
$wsQ1.Range("D3").Value = "7.93"
$wsQ1.Range("E3").Value = "62.03"
$wsQ1.Range("F3").Value = "1.74"
$wsQ1.Range("G3").Value = "0.1380"
$wsQ1.Range("H3").Value = 9

$wsQ1.Range("B4").Value = "009055"
$wsQ1.Range("C4").Value = "圆信永丰大湾区主题混合A"
$wsQ1.Range("D4").Value = "1.60"
$wsQ1.Range("E4").Value = "92.94"
$wsQ1.Range("F4").Value = "2.70"
$wsQ1.Range("G4").Value = "0.0432"
$wsQ1.Range("H4").Value = 10

$wsQ1.Range("B5").Value = "004265"
$wsQ1.Range("C5").Value = "金鹰民丰回报定期开放混合"
$wsQ1.Range("D5").Value = "6.57"
$wsQ1.Range("E5").Value = "28.61"
$wsQ1.Range("F5").Value = "0.64"
$wsQ1.Range("G5").Value = "0.0420"
$wsQ1.Range("H5").Value = 10

$wsQ1.Range("B6").Value = "009568"
$wsQ1.Range("C6").Value = "浙商智多宝稳健一年持有期混合A"
$wsQ1.Range("D6").Value = "3.13"
$wsQ1.Range("E6").Value = "22.64"
$wsQ1.Range("F6").Value = "1.28"
$wsQ1.Range("G6").Value = "0.0401"
$wsQ1.Range("H6").Value = 5

$wsQ1.Range("B7").Value = "009056"
$wsQ1.Range("C7").Value = "圆信永丰大湾区主题混合C"
$wsQ1.Range("D7").Value = "1.24"
$wsQ1.Range("E7").Value = "92.94"
$wsQ1.Range("F7").Value = "2.70"
$wsQ1.Range("G7").Value = "0.0335"
$wsQ1.Range("H7").Value = 10

$wsQ1.Range("B8").Value = "009569"
$wsQ1.Range("C8").Value = "浙商智多宝稳健一年持有期混合C"
$wsQ1.Range("D8").Value = "1.59"
$wsQ1.Range("E8").Value = "22.64"
$wsQ1.Range("F8").Value = "1.28"
$wsQ1.Range("G8").Value = "0.0204"
$wsQ1.Range("H8").Value = 5

# ---------------------------------------------------------------------
# 2. Re-create the "总计" (summary) sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

$wsQ4.Range("B1:D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)

$wsQ4.Range("A2").Copy()
$wsTotal.Range("A2:A3").PasteSpecial(-4122)

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 0.79

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.06

# Restore the original active tab ("2021-Q4", first sheet).
$wsQ4.Activate()
